$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (reddit_50k)
$ws.Range("G2").Value = 10.83716549873352
$ws.Range("H2").Value = 62.82888569831848
$ws.Range("I2").Value = 3.744938850402832
$ws.Range("J2").Value = 77.67947835922241

# Row 3 (reddit_100k)
$ws.Range("G3").Value = 9.58975796699524
$ws.Range("H3").Value = 28.7363639831543
$ws.Range("I3").Value = 2.95425181388855
$ws.Range("J3").Value = 41.57203869819641

# Row 4 (reddit_200k)
$ws.Range("G4").Value = 13.54277448654175
$ws.Range("H4").Value = 46.2468533039093
$ws.Range("I4").Value = 4.371213436126709
$ws.Range("J4").Value = 64.38922328948975

# Row 5 (reddit_500k)
$ws.Range("G5").Value = 29.63773217201233
$ws.Range("H5").Value = 232.3113450527191
$ws.Range("I5").Value = 24.56981449127197
$ws.Range("J5").Value = 286.7209562778473

# Row 6 (corpus-webis-tldr-17)
$ws.Range("G6").Value = 104.0061954975128
$ws.Range("H6").Value = 618.1843870639801
$ws.Range("I6").Value = 84.07183785438538
$ws.Range("J6").Value = 806.6028195858001
